# Fixed naive component forecaster bug - Presentation state 11.02.
# A new first error-column ("AVERAGE_1") was inserted before the existing
# AVERAGE_1..AVERAGE_9 data on every row (r=2..24). Every existing value in a
# row shifts one column to the right (B->C, C->D, ... J->K); for rows that
# already used all 10 data columns (B:K) the previously last value (old K)
# is dropped since the sheet has no column beyond K for this table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.765134014560337
$ws.Cells.Item(2, 3).Value = 9.634600668846847
$ws.Cells.Item(2, 4).Value = -8.622337597701534
$ws.Cells.Item(2, 5).Value = -0.8522339323542271
$ws.Cells.Item(2, 6).Value = 0.478975912873543
$ws.Cells.Item(2, 7).Value = -1.892429927382574
$ws.Cells.Item(2, 8).Value = -1.228164494743756
$ws.Cells.Item(2, 9).Value = -0.4586921403196634
$ws.Cells.Item(2, 10).Value = -0.5513017252472258
$ws.Cells.Item(2, 11).Value = 0.1696252062987764

$ws.Cells.Item(3, 2).Value = 7.384797829235454
$ws.Cells.Item(3, 3).Value = -10.87214043731293
$ws.Cells.Item(3, 4).Value = -3.102036771965619
$ws.Cells.Item(3, 5).Value = -1.770826926737849
$ws.Cells.Item(3, 6).Value = -4.142232766993966
$ws.Cells.Item(3, 7).Value = -3.477967334355148
$ws.Cells.Item(3, 8).Value = -2.708494979931056
$ws.Cells.Item(3, 9).Value = -2.801104564858618
$ws.Cells.Item(3, 10).Value = -2.080177633312616
$ws.Cells.Item(3, 11).Value = -2.78912955792647

$ws.Cells.Item(4, 2).Value = -18.34129076817022
$ws.Cells.Item(4, 3).Value = -10.57118710282291
$ws.Cells.Item(4, 4).Value = -9.239977257595143
$ws.Cells.Item(4, 5).Value = -11.61138309785126
$ws.Cells.Item(4, 6).Value = -10.94711766521244
$ws.Cells.Item(4, 7).Value = -10.17764531078835
$ws.Cells.Item(4, 8).Value = -10.27025489571591
$ws.Cells.Item(4, 9).Value = -9.549327964169908
$ws.Cells.Item(4, 10).Value = -10.25827988878376
$ws.Cells.Item(4, 11).Value = -9.766173896461286

$ws.Cells.Item(5, 2).Value = 7.600074697314557
$ws.Cells.Item(5, 3).Value = 8.931284542542326
$ws.Cells.Item(5, 4).Value = 6.55987870228621
$ws.Cells.Item(5, 5).Value = 7.224144134925027
$ws.Cells.Item(5, 6).Value = 7.99361648934912
$ws.Cells.Item(5, 7).Value = 7.901006904421558
$ws.Cells.Item(5, 8).Value = 8.62193383596756
$ws.Cells.Item(5, 9).Value = 7.912981911353705
$ws.Cells.Item(5, 10).Value = 8.405087903676183
$ws.Cells.Item(5, 11).Value = 8.87489265094917

$ws.Cells.Item(6, 2).Value = 0.7916129955631771
$ws.Cells.Item(6, 3).Value = -1.579792844692939
$ws.Cells.Item(6, 4).Value = -0.9155274120541215
$ws.Cells.Item(6, 5).Value = -0.1460550576300293
$ws.Cells.Item(6, 6).Value = -0.2386646425575917
$ws.Cells.Item(6, 7).Value = 0.4822622889884105
$ws.Cells.Item(6, 8).Value = -0.2266896356254442
$ws.Cells.Item(6, 9).Value = 0.2654163566970329
$ws.Cells.Item(6, 10).Value = 0.73522110397002
$ws.Cells.Item(6, 11).Value = 1.048046298935466

$ws.Cells.Item(7, 2).Value = -3.727363316492332
$ws.Cells.Item(7, 3).Value = -3.063097883853514
$ws.Cells.Item(7, 4).Value = -2.293625529429422
$ws.Cells.Item(7, 5).Value = -2.386235114356984
$ws.Cells.Item(7, 6).Value = -1.665308182810982
$ws.Cells.Item(7, 7).Value = -2.374260107424837
$ws.Cells.Item(7, 8).Value = -1.88215411510236
$ws.Cells.Item(7, 9).Value = -1.412349367829372
$ws.Cells.Item(7, 10).Value = -1.099524172863926
$ws.Cells.Item(7, 11).Value = -1.977926754115506

$ws.Cells.Item(8, 2).Value = 0.376932102669816
$ws.Cells.Item(8, 3).Value = 1.146404457093908
$ws.Cells.Item(8, 4).Value = 1.053794872166346
$ws.Cells.Item(8, 5).Value = 1.774721803712348
$ws.Cells.Item(8, 6).Value = 1.065769879098493
$ws.Cells.Item(8, 7).Value = 1.55787587142097
$ws.Cells.Item(8, 8).Value = 2.027680618693958
$ws.Cells.Item(8, 9).Value = 2.340505813659403
$ws.Cells.Item(8, 10).Value = 1.462103232407824
$ws.Cells.Item(8, 11).Value = 0.7321912172293545

$ws.Cells.Item(9, 2).Value = 1.207578635508109
$ws.Cells.Item(9, 3).Value = 1.114969050580547
$ws.Cells.Item(9, 4).Value = 1.835895982126549
$ws.Cells.Item(9, 5).Value = 1.126944057512694
$ws.Cells.Item(9, 6).Value = 1.619050049835171
$ws.Cells.Item(9, 7).Value = 2.088854797108159
$ws.Cells.Item(9, 8).Value = 2.401679992073604
$ws.Cells.Item(9, 9).Value = 1.523277410822025
$ws.Cells.Item(9, 10).Value = 0.7933653956435556
$ws.Cells.Item(9, 11).Value = 2.455544492033183

$ws.Cells.Item(10, 2).Value = -0.9264868865757077
$ws.Cells.Item(10, 3).Value = -0.2055599550297054
$ws.Cells.Item(10, 4).Value = -0.9145118796435601
$ws.Cells.Item(10, 5).Value = -0.422405887321083
$ws.Cells.Item(10, 6).Value = 0.04739885995190407
$ws.Cells.Item(10, 7).Value = 0.36022405491735
$ws.Cells.Item(10, 8).Value = -0.5181785263342299
$ws.Cells.Item(10, 9).Value = -1.248090541512699
$ws.Cells.Item(10, 10).Value = 0.4140885548769285
$ws.Cells.Item(10, 11).Value = -0.2051307335183153

$ws.Cells.Item(11, 2).Value = 0.3770345820039356
$ws.Cells.Item(11, 3).Value = -0.3319173426099191
$ws.Cells.Item(11, 4).Value = 0.160188649712558
$ws.Cells.Item(11, 5).Value = 0.6299933969855451
$ws.Cells.Item(11, 6).Value = 0.942818591950991
$ws.Cells.Item(11, 7).Value = 0.06441601069941108
$ws.Cells.Item(11, 8).Value = -0.6654960044790579
$ws.Cells.Item(11, 9).Value = 0.9966830919105695
$ws.Cells.Item(11, 10).Value = 0.3774638035153257
$ws.Cells.Item(11, 11).Value = 0.6147675671350392

$ws.Cells.Item(12, 2).Value = -0.4275923834192769
$ws.Cells.Item(12, 3).Value = 0.0645136089032002
$ws.Cells.Item(12, 4).Value = 0.5343183561761873
$ws.Cells.Item(12, 5).Value = 0.8471435511416332
$ws.Cells.Item(12, 6).Value = -0.03125903010994671
$ws.Cells.Item(12, 7).Value = -0.7611710452884157
$ws.Cells.Item(12, 8).Value = 0.9010080511012117
$ws.Cells.Item(12, 9).Value = 0.2817887627059679
$ws.Cells.Item(12, 10).Value = 0.5190925263256815
$ws.Cells.Item(12, 11).Value = 0.4158151564502698

$ws.Cells.Item(13, 2).Value = 0.324932645901923
$ws.Cells.Item(13, 3).Value = 0.7947373931749101
$ws.Cells.Item(13, 4).Value = 1.107562588140356
$ws.Cells.Item(13, 5).Value = 0.2291600068887761
$ws.Cells.Item(13, 6).Value = -0.5007520082896928
$ws.Cells.Item(13, 7).Value = 1.161427088099934
$ws.Cells.Item(13, 8).Value = 0.5422077997046907
$ws.Cells.Item(13, 9).Value = 0.7795115633244043
$ws.Cells.Item(13, 10).Value = 0.6762341934489926
$ws.Cells.Item(13, 11).Value = -0.1153642338804421

$ws.Cells.Item(14, 2).Value = -0.04071760298358112
$ws.Cells.Item(14, 3).Value = 0.2721075919818648
$ws.Cells.Item(14, 4).Value = -0.6062949892697151
$ws.Cells.Item(14, 5).Value = -1.336207004448184
$ws.Cells.Item(14, 6).Value = 0.3259720919414433
$ws.Cells.Item(14, 7).Value = -0.2932471964538005
$ws.Cells.Item(14, 8).Value = -0.05594343283408693
$ws.Cells.Item(14, 9).Value = -0.1592208027094986
$ws.Cells.Item(14, 10).Value = -0.9508192300389333
$ws.Cells.Item(14, 11).Value = -0.3403303223714723

$ws.Cells.Item(15, 2).Value = 0.3721869518844864
$ws.Cells.Item(15, 3).Value = -0.5062156293670936
$ws.Cells.Item(15, 4).Value = -1.236127644545562
$ws.Cells.Item(15, 5).Value = 0.4260514518440648
$ws.Cells.Item(15, 6).Value = -0.193167836551179
$ws.Cells.Item(15, 7).Value = 0.04413592706853459
$ws.Cells.Item(15, 8).Value = -0.05914144280687711
$ws.Cells.Item(15, 9).Value = -0.8507398701363118
$ws.Cells.Item(15, 10).Value = -0.2402509624688508
$ws.Cells.Item(15, 11).Value = -0.4017729932881683

$ws.Cells.Item(16, 2).Value = -0.1524291232873974
$ws.Cells.Item(16, 3).Value = -0.8823411384658664
$ws.Cells.Item(16, 4).Value = 0.779837957923761
$ws.Cells.Item(16, 5).Value = 0.1606186695285172
$ws.Cells.Item(16, 6).Value = 0.3979224331482308
$ws.Cells.Item(16, 7).Value = 0.2946450632728191
$ws.Cells.Item(16, 8).Value = -0.4969533640566156
$ws.Cells.Item(16, 9).Value = 0.1135355436108454
$ws.Cells.Item(16, 10).Value = -0.04798648720847212

$ws.Cells.Item(17, 2).Value = -1.030518528898312
$ws.Cells.Item(17, 3).Value = 0.6316605674913157
$ws.Cells.Item(17, 4).Value = 0.0124412790960719
$ws.Cells.Item(17, 5).Value = 0.2497450427157855
$ws.Cells.Item(17, 6).Value = 0.1464676728403738
$ws.Cells.Item(17, 7).Value = -0.6451307544890609
$ws.Cells.Item(17, 8).Value = -0.03464184682159993
$ws.Cells.Item(17, 9).Value = -0.1961638776409175

$ws.Cells.Item(18, 2).Value = 0.4742145784871607
$ws.Cells.Item(18, 3).Value = -0.1450047099080831
$ws.Cells.Item(18, 4).Value = 0.0922990537116305
$ws.Cells.Item(18, 5).Value = -0.0109783161637812
$ws.Cells.Item(18, 6).Value = -0.8025767434932158
$ws.Cells.Item(18, 7).Value = -0.1920878358257549
$ws.Cells.Item(18, 8).Value = -0.3536098666450724

$ws.Cells.Item(19, 2).Value = 0.3556547466179877
$ws.Cells.Item(19, 3).Value = 0.5929585102377013
$ws.Cells.Item(19, 4).Value = 0.4896811403622896
$ws.Cells.Item(19, 5).Value = -0.3019172869671451
$ws.Cells.Item(19, 6).Value = 0.3085716207003159
$ws.Cells.Item(19, 7).Value = 0.1470495898809984

$ws.Cells.Item(20, 2).Value = 0.3126006297022321
$ws.Cells.Item(20, 3).Value = 0.2093232598268204
$ws.Cells.Item(20, 4).Value = -0.5822751675026142
$ws.Cells.Item(20, 5).Value = 0.02821374016484672
$ws.Cells.Item(20, 6).Value = -0.1333082906544708

$ws.Cells.Item(21, 2).Value = 0.3812981176718321
$ws.Cells.Item(21, 3).Value = -0.4103003096576026
$ws.Cells.Item(21, 4).Value = 0.2001885980098584
$ws.Cells.Item(21, 5).Value = 0.03866656719054083

$ws.Cells.Item(22, 2).Value = -0.716162849403934
$ws.Cells.Item(22, 3).Value = -0.1056739417364731
$ws.Cells.Item(22, 4).Value = -0.2671959725557906

$ws.Cells.Item(23, 2).Value = 0.506656010950813
$ws.Cells.Item(23, 3).Value = 0.3451339801314955

$ws.Cells.Item(24, 2).Value = -0.343237405067616
